# "MVP Statistics.xlsx" - TeamMVP sheet update
#
# Summary of the change being applied (per the commit's OOXML diff):
#   1. Header cell C1 is renamed from "Count" to "Winner Count".
#   2. The Los Angeles Angels row (row 18) has its franchise name trimmed
#      from "Los Angeles Angels of Anaheim " to "Los Angeles Angels ".
#   3. A brand-new column D, "Top 9 Contenders", is added with a per-team
#      count value for every data row (2-38).
#
# All other existing data (columns A/B "Team"/franchise name and C "Winner
# Count" numbers) is left exactly as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TeamMVP")

# --- 1. Rename the existing "Count" header to "Winner Count" -------------
$ws.Range("C1").Value = "Winner Count"

# --- 2. Add the new "Top 9 Contenders" column (D) -------------------------
$ws.Range("D1").Value = "Top 9 Contenders"

# --- 3. Fix the LA Angels franchise name (row 18) -------------------------
$ws.Range("B18").Value = "Los Angeles Angels "

$top9 = @{
    2  = 4
    3  = 39
    4  = 8
    5  = 41
    6  = 70
    7  = 3
    8  = 9
    9  = 32
    10 = 35
    11 = 52
    12 = 34
    13 = 21
    14 = 38
    15 = 5
    16 = 33
    17 = 22
    18 = 18
    19 = 53
    20 = 2
    21 = 26
    22 = 36
    23 = 11
    24 = 19
    25 = 1
    26 = 28
    27 = 69
    28 = 29
    29 = 36
    30 = 39
    31 = 21
    32 = 21
    33 = 46
    34 = 54
    35 = 7
    36 = 23
    37 = 29
    38 = 17
}

foreach ($row in $top9.Keys | Sort-Object) {
    $ws.Range("D$row").Value = $top9[$row]
}
